# Update LUTs for more precise binning borders.
#
# N5 is the anchor of a shared formula group that spans N5:N43
# (N5 itself, then N6 as the shared "master" with ref N6:N43, si=8).
# The edit subtracts 1 from the J-column value before it gets
# concatenated into the generated C-array text, e.g.
#   J5&", "&...   ->   J5-1&", "&...
#
# Re-assigning FormulaR1C1 across the whole N5:N43 range in one shot
# keeps Excel's shared-formula grouping intact (master cell + "t=shared"
# followers) instead of turning every row into its own literal formula.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("N5:N43").FormulaR1C1 = '=IF(RC[-4]<>"",RC[-4]-1&", "&IF(RC8<>R[1]C8,CHAR(10),""),256^R8C3-1&CHAR(10))'

# Recalculate so every dependent cell (P5, B16, ...) picks up the new
# LUT values.
$excel.CalculateFull()

# Reflect the author's final viewport/selection: scrolled down a bit and
# with S26 selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("S26").Select()
